# Auto-generated Word COM-interop script.
# Applies a series of literal text find & replace edits to $word.ActiveDocument,
# mirroring the author's copy-edit / grammar-fix commit.

$d = $word.ActiveDocument
$appliedCount = 0
$failedFinds = @()

function Apply-Replace($FindText, $ReplaceText, $ReplaceMode) {
    $range = $d.Content
    $ok = $range.Find.Execute(
        $FindText,    # FindText
        $false,       # MatchCase
        $false,       # MatchWholeWord
        $false,       # MatchWildcards
        $false,       # MatchSoundsLike
        $false,       # MatchAllWordForms
        $true,        # Forward
        1,            # Wrap (wdFindContinue)
        $false,       # Format
        $ReplaceText, # ReplaceWith
        $ReplaceMode  # Replace (1 = wdReplaceOne, 2 = wdReplaceAll)
    )

    if ($ok) {
        $appliedCount = $appliedCount + 1
    } else {
        $failedFinds += $FindText
    }
}

# intro client-side
Apply-Replace " within the present work. The client side of " " within the present work. The client-side of " 1

# add 'with'
Apply-Replace "called Declarative API. Its main purpose is to provide the user an asynchronous interface" "called Declarative API. Its main purpose is to provide the user with an asynchronous interface" 1

# case study which->that
Apply-Replace "together with a case study which involves the implementation of" "together with a case study that involves the implementation of" 1

# and in addition -> , and also
Apply-Replace "data analysis, data transfers, data management and, in addition, features like staging data from tape" "data analysis, data transfers, data management, and also features like staging data from tape" 1

# remove comma after Although
Apply-Replace "sites that are hosting the data. Although, XRootD supported multi-storage deployments" "sites that are hosting the data. Although XRootD supported multi-storage deployments" 1

# feature which->that
Apply-Replace "the addition of a feature which allowed its proper functionality within a global, multi-site environment was in fact the core idea of AAA." "the addition of a feature that allowed its proper functionality within a global, multi-site environment was in fact the core idea of AAA." 1

# In order to -> To / add the
Apply-Replace "In order to emphasize the importance of XRootD, it is worth mentioning" "To emphasize the importance of the XRootD, it is worth mentioning" 1

# server side/client side hyphenate + the WLCG
Apply-Replace " framework, both in terms of its server side as well as its client side, since both implementations are crucial in understanding the overall workflow of data access and data manipulation within WLCG community." " framework, both in terms of its server-side as well as its client-side, since both implementations are crucial in understanding the overall workflow of data access and data manipulation within the WLCG community." 1

# experiments which->that ran
Apply-Replace "The main objective of any scientific project that is based on experiments which ran at CERN" "The main objective of any scientific project that is based on experiments that ran at CERN" 1

# compute->computing resources
Apply-Replace "is the access to the compute resources which are used for submitting jobs" "is the access to the computing resources which are used for submitting jobs" 1

# diagram below -> Fig.1
Apply-Replace "An old model of such a workflow is described in diagram below (also called “jobs go to data” paradigm [2])." "An old model of such a workflow is described in Fig.1 (also called “jobs go to data” paradigm [2])." 1

# add space + comma result,
Apply-Replace "the grid middleware [2]).This introduces significant overhead and as a result slows down the process" "the grid middleware [2]). This introduces significant overhead and as a result, slows down the process" 1

# Federated -> A Federated ; which are cooperating with each other -> that are cooperating
Apply-Replace "with the required data. Federated storage system is the implementation that aims at solving such issues. Defined in [2] as a collection of unpaired storage resources that are managed by a set of domains which are cooperating with each other (but also independent) and also are accessible via a common namespace. " "with the required data. A Federated storage system is the implementation that aims at solving such issues. Defined in [2] as a collection of unpaired storage resources that are managed by a set of domains that are cooperating (but also independent) and also are accessible via a common namespace. " 1

# remove 'a' direct contact
Apply-Replace "makes a direct contact with the central endpoint" "makes direct contact with the central endpoint" 1

# cache mechanism, and -> A tremendous / recent years
Apply-Replace "cache mechanism and many more. In fact, a tremendous work (progress) has been done in the recent years, especially for extending the scalability features" "cache mechanism, and many more. A tremendous work (progress) has been done in recent years, especially for extending the scalability features" 1

# [2], and / The development team
Apply-Replace " [2] and many more. In fact, the development team is constantly committing new or improved " " [2], and many more. The development team is constantly committing new or improved " 1

# Oxford comma layers
Apply-Replace ", namely: Network layer, Protocol layer, File-system layer and a Storage laye" ", namely: Network layer, Protocol layer, File-system layer, and a Storage laye" 1

# added into -> added to
Apply-Replace "run-time plug-in mechanism, new features can be added into the framework with little effort." "run-time plug-in mechanism, new features can be added to the framework with little effort." 1

# remove comma after implementation
Apply-Replace "The XrdCl implementation, developed within a multi-threaded " "The XrdCl implementation developed within a multi-threaded " 1

# add 'a' slight increase
Apply-Replace "(however, at the cost of slight increase in code complexity)." "(however, at the cost of a slight increase in code complexity)." 1

# done on -> done in a single-threaded manner
Apply-Replace "This execution of requests is done on a single-threaded manner" "This execution of requests is done in a single-threaded manner" 1

# add 'the' time-to-live
Apply-Replace "until time-to-live timeout elapses." "until the time-to-live timeout elapses." 1

# add 'a' network of WAN
Apply-Replace "improving in this way the performance over network of WAN type." "improving in this way the performance over a network of WAN type." 1

# reestablish -> re-establish
Apply-Replace "for example, the user might want to reestablish the connection with new credentials." "for example, the user might want to re-establish the connection with new credentials." 1

# are organized -> is organized
Apply-Replace ", where each of the components are organized in three main categories: " ", where each of the components is organized in three main categories: " 1

# Oxford comma External
Apply-Replace "XRootD-Core, XRootD and External." "XRootD-Core, XRootD, and External." 1

# add 'the' XRootD client
Apply-Replace "The XrdCl library is the foundation part of the following components of XRootD client:" "The XrdCl library is the foundation part of the following components of the XRootD client:" 1

# command line -> command-line
Apply-Replace "The command line interface [6]." "The command-line interface [6]." 1

# typo fix Asyncrhonous->Asynchronous
Apply-Replace "Asyncrhonous implementations within XRootD client interface" "Asynchronous implementations within XRootD client interface" 1

# remove clause about file access API
Apply-Replace " is mainly used with file-based data repositories, a crucial component is indeed the file access API, that contains both single file as well as file system implementations. It was already mentioned that these objects have both synchronous and asynchronous behavior. " " is mainly used with file-based data repositories, a crucial component is indeed the file access API. It was already mentioned that these objects have both synchronous and asynchronous behavior. " 1

# flow which->that involves
Apply-Replace "A usual execution flow which involves file operations" "A usual execution flow that involves file operations" 1

# the existing API / end-user
Apply-Replace " It has been built on top of existing API and provides an additional layer of abstraction (that layer itself is what makes a more convenient interface between the client and the end user)." " It has been built on top of the existing API and provides an additional layer of abstraction (that layer itself is what makes a more convenient interface between the client and the end-user)." 1

# remove 'actually'
Apply-Replace "Its key features that actually make the API easy to use are the following:" "Its key features that make the API easy to use are the following:" 1

# Syntax -> The syntax / on->to execution flow
Apply-Replace "Syntax is declarative-centric, meaning that users should focus on the actual choice of operation rather than paying much attention (effort) on the execution flow." "The syntax is declarative-centric, meaning that users should focus on the actual choice of operation rather than paying much attention (effort) to the execution flow." 1

# add 'the' user / 'the' compilation phase
Apply-Replace "Proper signaling for user of any incorrect declarations and configurations during compilation phase." "Proper signaling for the user of any incorrect declarations and configurations during the compilation phase." 1

# done in a consistent manner -> done consistently
Apply-Replace "Error handling for the workflow is done in a consistent manner, showing proper error messages" "Error handling for the workflow is done consistently, showing proper error messages" 1

# add 'the' result of one operation
Apply-Replace "The constructed API makes it so there is a communication protocol between the operations: result of one operation is used to compute the following operation, making this implementation very robust." "The constructed API makes it so there is a communication protocol between the operations: the result of one operation is used to compute the following operation, making this implementation very robust." 1

# add 'the' modern C++ language paradigm
Apply-Replace " One can see that the new API is more in line with modern C++ language paradigm." " One can see that the new API is more in line with the modern C++ language paradigm." 1

# Oxford comma Open Read and Close
Apply-Replace "Open, Read and Close) the corresponding object is created and the operator () is used" "Open, Read, and Close) the corresponding object is created and the operator () is used" 1

# add 'the' last line
Apply-Replace "; last line contains a utility for synchronous execution of the pipeline (current " "; the last line contains a utility for synchronous execution of the pipeline (current " 1

# connected between->to each other
Apply-Replace "The defined operations are connected between each other by the | operator. " "The defined operations are connected to each other by the | operator. " 1

# add 'the' C++ programming language
Apply-Replace "because of the operator overloading feature of C++ programming language." "because of the operator overloading feature of the C++ programming language." 1

# remove 'a' parallel execution
Apply-Replace "The syntax also supports a parallel execution of multiple flows of operations" "The syntax also supports parallel execution of multiple flows of operations" 1

# Example -> An example pipeline
Apply-Replace "It also accepts a variable number of operations. Example pipeline with three parallel operations can be seen in Listing 6." "It also accepts a variable number of operations. An example pipeline with three parallel operations can be seen in Listing 6." 1

# user wants -> the user wants
Apply-Replace "the following example is proposed: user wants to access a file with a size of 0.5MB from a data batch" "the following example is proposed: the user wants to access a file with a size of 0.5MB from a data batch" 1

# remove 'the' Listing 8
Apply-Replace ". Using the declarative approach, the procedure will look like the Listing 8." ". Using the declarative approach, the procedure will look like Listing 8." 1

# add 'the' declaration of a lock file
Apply-Replace " The first line is for declaration of a lock file," " The first line is for the declaration of a lock file," 1

# Oxford comma open a read and a close
Apply-Replace "the pipeline continues by doing an open, a read and a close for the actual" "the pipeline continues by doing an open, a read, and a close for the actual" 1

# remove comma before since it is not needed
Apply-Replace "file that needs to be accessed. The Rm function is used for deleting the lock file, since it is not needed anymore." "file that needs to be accessed. The Rm function is used for deleting the lock file since it is not needed anymore." 1

# remove 'an' error status
Apply-Replace "but their handlers will be called (with an error status). Using the pipelining semantic" "but their handlers will be called (with error status). Using the pipelining semantic" 1

# add comma before 'and stored'
Apply-Replace "expanded and encoded with redundant data pieces and stored across a set of different locations or storage media." "expanded and encoded with redundant data pieces, and stored across a set of different locations or storage media." 1

# striped -> stripped
Apply-Replace "will imply that an entire block of data will be striped into " "will imply that an entire block of data will be stripped into " 1

# Figure -> Fig. (caption, 2 occurrences)
Apply-Replace "Figure 4: Asynchronous workflow in terms of response handlers." "Fig. 4: Asynchronous workflow in terms of response handlers." 2

# on parallel -> in parallel
Apply-Replace "be done on parallel. Listing 9 contains the necessary workflow" "be done in parallel. Listing 9 contains the necessary workflow" 1

# server side/client side -> hyphenated (conclusion)
Apply-Replace "A short description of the architecture for both the server side as well as the client side was discussed." "A short description of the architecture for both the server-side as well as the client-side was discussed." 1

# add 'an' asynchronous / remove comma before while
Apply-Replace "It is showed that Declarative API is an efficient tool in providing asynchronous C++ interface for the user, while keeping a clear and concise workflow." "It is showed that Declarative API is an efficient tool in providing an asynchronous C++ interface for the user while keeping a clear and concise workflow." 1

# add 'the' department
Apply-Replace " – head of department, and " " – head of the department, and " 1


Write-Output "Applied: $appliedCount / 57"
if ($failedFinds.Count -gt 0) {
    Write-Output "FAILED FINDS:"
    foreach ($f in $failedFinds) {
        Write-Output "  - $f"
    }
}
